# Amend corrected label annotations: normalize casing / ordering of the
# "labels" (column F) values on sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = "135_product_information"
    9   = "135_product_information"
    11  = "application instructions"
    12  = "use restrictions || application instructions"
    13  = "use restrictions || application instructions"
    14  = "use restrictions || application instructions"
    15  = "use restrictions || application instructions"
    16  = "use restrictions || application instructions"
    17  = "use restrictions || application instructions"
    18  = "use restrictions || application instructions"
    19  = "use restrictions"
    20  = "use restrictions"
    21  = "use restrictions || application instructions"
    22  = "use restrictions || application instructions"
    29  = "93_referral_statement"
    38  = "18_hazards_to_humans_and_domestic_animals"
    39  = "ppe"
    40  = "ppe"
    45  = "off target movement || application instructions || env warning - species || env warning - water"
    52  = "application instructions"
    53  = "135_product_information"
    54  = "application instructions"
    59  = "application instructions"
    60  = "application instructions"
    61  = "application instructions"
    62  = "mixing || application instructions"
    63  = "application instructions"
    64  = "irrigation || application instructions || chemigation"
    65  = "off target movement"
    66  = "irrigation"
    67  = "irrigation || chemigation"
    69  = "use restrictions"
    70  = "use restrictions"
    249 = "application instructions"
    250 = "mixing"
    251 = "mixing"
    252 = "mixing"
    261 = "154_pesticide_storage"
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
